$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 20.65176033333333
$ws.Range("H2").Value = 61.955281
$ws.Range("I2").Value = 0.8502439951095444
$ws.Range("J2").Value = 0.8502439951095444
$ws.Range("M2").Value = 1.910418
$ws.Range("N2").Value = 5.731254
$ws.Range("O2").Value = 0.01809124304049503
$ws.Range("P2").Value = 0.01809124304049503
$ws.Range("Q2").Value = 39.45349467248599
$ws.Range("R2").Value = 355.081452052374
$ws.Range("S2").Value = 0.01538197075924823
$ws.Range("T2").Value = 0.01538197075924823

$ws.Range("G3").Value = 20.65176033333333
$ws.Range("H3").Value = 61.955281
$ws.Range("I3").Value = 0.8502439951095444
$ws.Range("J3").Value = 0.8502439951095444
$ws.Range("O3").Value = 0.302988173785169
$ws.Range("P3").Value = 0.302988173785169
$ws.Range("Q3").Value = 660.7584826262073
$ws.Range("R3").Value = 5946.826343635867
$ws.Range("S3").Value = 0.257613875350047
$ws.Range("T3").Value = 0.257613875350047

$ws.Range("G4").Value = 20.65176033333333
$ws.Range("H4").Value = 61.955281
$ws.Range("I4").Value = 0.8502439951095444
$ws.Range("J4").Value = 0.8502439951095444
$ws.Range("M4").Value = 37.858701
$ws.Range("N4").Value = 113.576103
$ws.Range("O4").Value = 0.3585136661130873
$ws.Range("P4").Value = 0.3585136661130873
$ws.Range("Q4").Value = 781.8488195833269
$ws.Range("R4").Value = 7036.639376249942
$ws.Range("S4").Value = 0.3048240917773606
$ws.Range("T4").Value = 0.3048240917773606

$ws.Range("G5").Value = 20.65176033333333
$ws.Range("H5").Value = 61.955281
$ws.Range("I5").Value = 0.8502439951095444
$ws.Range("J5").Value = 0.8502439951095444
$ws.Range("M5").Value = 33.83466466666667
$ws.Range("N5").Value = 101.503994
$ws.Range("O5").Value = 0.3204069170612486
$ws.Range("P5").Value = 0.3204069170612486
$ws.Range("Q5").Value = 698.7453856547015
$ws.Range("R5").Value = 6288.708470892314
$ws.Range("S5").Value = 0.2724240572228884
$ws.Range("T5").Value = 0.2724240572228884

$ws.Range("G6").Value = 0.3109896666666667
$ws.Range("H6").Value = 0.932969
$ws.Range("I6").Value = 0.01280361055699766
$ws.Range("J6").Value = 0.01280361055699766
$ws.Range("M6").Value = 1.910418
$ws.Range("N6").Value = 5.731254
$ws.Range("O6").Value = 0.01809124304049503
$ws.Range("P6").Value = 0.01809124304049503
$ws.Range("Q6").Value = 0.594120257014
$ws.Range("R6").Value = 5.347082313126
$ws.Range("S6").Value = 0.0002316332303824926
$ws.Range("T6").Value = 0.0002316332303824926

$ws.Range("G7").Value = 0.3109896666666667
$ws.Range("H7").Value = 0.932969
$ws.Range("I7").Value = 0.01280361055699766
$ws.Range("J7").Value = 0.01280361055699766
$ws.Range("O7").Value = 0.302988173785169
$ws.Range("P7").Value = 0.302988173785169
$ws.Range("Q7").Value = 9.950195864292667
$ws.Range("R7").Value = 89.55176277863401
$ws.Range("S7").Value = 0.003879342580521232
$ws.Range("T7").Value = 0.003879342580521232

$ws.Range("G8").Value = 0.3109896666666667
$ws.Range("H8").Value = 0.932969
$ws.Range("I8").Value = 0.01280361055699766
$ws.Range("J8").Value = 0.01280361055699766
$ws.Range("M8").Value = 37.858701
$ws.Range("N8").Value = 113.576103
$ws.Range("O8").Value = 0.3585136661130873
$ws.Range("P8").Value = 0.3585136661130873
$ws.Range("Q8").Value = 11.773664804423
$ws.Range("R8").Value = 105.962983239807
$ws.Range("S8").Value = 0.00459026936027346
$ws.Range("T8").Value = 0.00459026936027346

$ws.Range("G9").Value = 0.3109896666666667
$ws.Range("H9").Value = 0.932969
$ws.Range("I9").Value = 0.01280361055699766
$ws.Range("J9").Value = 0.01280361055699766
$ws.Range("M9").Value = 33.83466466666667
$ws.Range("N9").Value = 101.503994
$ws.Range("O9").Value = 0.3204069170612486
$ws.Range("P9").Value = 0.3204069170612486
$ws.Range("Q9").Value = 10.52223108646511
$ws.Range("R9").Value = 94.70007977818601
$ws.Range("S9").Value = 0.004102365385820476
$ws.Range("T9").Value = 0.004102365385820476

$ws.Range("G10").Value = 3.326466333333334
$ws.Range("H10").Value = 9.979399000000001
$ws.Range("I10").Value = 0.1369523943334579
$ws.Range("J10").Value = 0.1369523943334579
$ws.Range("M10").Value = 1.910418
$ws.Range("N10").Value = 5.731254
$ws.Range("O10").Value = 0.01809124304049503
$ws.Range("P10").Value = 0.01809124304049503
$ws.Range("Q10").Value = 6.354941159594
$ws.Range("R10").Value = 57.194470436346
$ws.Range("S10").Value = 0.002477639050864301
$ws.Range("T10").Value = 0.002477639050864301

$ws.Range("G11").Value = 3.326466333333334
$ws.Range("H11").Value = 9.979399000000001
$ws.Range("I11").Value = 0.1369523943334579
$ws.Range("J11").Value = 0.1369523943334579
$ws.Range("O11").Value = 0.302988173785169
$ws.Range("P11").Value = 0.302988173785169
$ws.Range("Q11").Value = 106.4311618691793
$ws.Range("R11").Value = 957.8804568226141
$ws.Range("S11").Value = 0.04149495585460075
$ws.Range("T11").Value = 0.04149495585460075

$ws.Range("G12").Value = 3.326466333333334
$ws.Range("H12").Value = 9.979399000000001
$ws.Range("I12").Value = 0.1369523943334579
$ws.Range("J12").Value = 0.1369523943334579
$ws.Range("M12").Value = 37.858701
$ws.Range("N12").Value = 113.576103
$ws.Range("O12").Value = 0.3585136661130873
$ws.Range("P12").Value = 0.3585136661130873
$ws.Range("Q12").Value = 125.935694300233
$ws.Range("R12").Value = 1133.421248702097
$ws.Range("S12").Value = 0.04909930497545321
$ws.Range("T12").Value = 0.04909930497545321

$ws.Range("G13").Value = 3.326466333333334
$ws.Range("H13").Value = 9.979399000000001
$ws.Range("I13").Value = 0.1369523943334579
$ws.Range("J13").Value = 0.1369523943334579
$ws.Range("M13").Value = 33.83466466666667
$ws.Range("N13").Value = 101.503994
$ws.Range("O13").Value = 0.3204069170612486
$ws.Range("P13").Value = 0.3204069170612486
$ws.Range("Q13").Value = 112.5498729132896
$ws.Range("R13").Value = 1012.948856219606
$ws.Range("S13").Value = 0.04388049445253966
$ws.Range("T13").Value = 0.04388049445253966
